$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the new, longer address text (stored width="29")
$ws.Columns.Item(1).ColumnWidth = 28.17

# Append new row 6 with the additional test data record.
# E6/G6/H6 look numeric but must stay text, so force text format first.
$ws.Range("E6").NumberFormat = "@"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("H6").NumberFormat = "@"

$ws.Range("A6").Value = "Adrress 25 f1 @#$%^&*!(#)#*"
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = "25 f1 first name"
$ws.Range("D6").Value = "25 f1 last name"
$ws.Range("E6").Value = "25012334567955"
$ws.Range("F6").Value = "25 F1 City"
$ws.Range("G6").Value = "251"
$ws.Range("H6").Value = "2501"
